# Diets by period.xlsx — update to "last versions of the files"
#  - rename sheets (BCS/wo-cereals naming -> OS/AS naming)
#  - replace shared-string food-item labels + add new items
#  - refresh gram/kcal figures per diet sheet (rows expand on several sheets)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet renames
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(1).Name = "I - OS"
$wb.Worksheets.Item(2).Name = "I - AS"
$wb.Worksheets.Item(3).Name = "II - OS"
$wb.Worksheets.Item(4).Name = "II - AS"
$wb.Worksheets.Item(5).Name = "III- OS"
$wb.Worksheets.Item(6).Name = "III - AS"

# ---------------------------------------------------------------------------
# 2) Sheet 1 ("I - OS")
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = 1417.0
$ws.Range("C2").Value = 2100.0

$ws.Range("A3").Value = "Potatoes"
$ws.Range("B3").Value = 344.8
$ws.Range("C3").Value = 300.0

$ws.Range("A4").Value = "Wheat flour"
$ws.Range("B4").Value = 68.7
$ws.Range("C4").Value = 250.0

$ws.Range("A5").Value = "Barley (pearled)"
$ws.Range("B5").Value = 24.4
$ws.Range("C5").Value = 30.0

$ws.Range("A6").Value = "Canola oil"
$ws.Range("B6").Value = 19.2
$ws.Range("C6").Value = 170.0

$ws.Range("A7").Value = "Rice (white)"
$ws.Range("B7").Value = 115.4
$ws.Range("C7").Value = 150.0

$ws.Range("A8").Value = "Corn flour (whole-grain)"
$ws.Range("B8").Value = 27.7
$ws.Range("C8").Value = 100.0

$ws.Range("A9").Value = "Corn"
$ws.Range("B9").Value = 104.2
$ws.Range("C9").Value = 100.0

$ws.Range("A10").Value = "Soy flour"
$ws.Range("B10").Value = 46.1
$ws.Range("C10").Value = 200.0

$ws.Range("A11").Value = "Soybeans"
$ws.Range("B11").Value = 116.3
$ws.Range("C11").Value = 200.0

$ws.Range("A12").Value = "Anchovy (raw)"
$ws.Range("B12").Value = 76.3
$ws.Range("C12").Value = 100.0

$ws.Range("A13").Value = "Cattle (lean)"
$ws.Range("B13").Value = 49.3
$ws.Range("C13").Value = 100.0

$ws.Range("A14").Value = "Cattle (organs)"
$ws.Range("B14").Value = 178.8
$ws.Range("C14").Value = 250.0

$ws.Range("A15").Value = "Milk (whole)"
$ws.Range("B15").Value = 245.9
$ws.Range("C15").Value = 150.0

# ---------------------------------------------------------------------------
# 3) Sheet 2 ("I - AS")
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = 1555.0
$ws.Range("C2").Value = 2100.0

$ws.Range("A3").Value = "Potatoes"
$ws.Range("B3").Value = 574.7
$ws.Range("C3").Value = 500.0

$ws.Range("A4").Value = "Anchovy (raw)"
$ws.Range("B4").Value = 114.5
$ws.Range("C4").Value = 150.0

$ws.Range("A5").Value = "Cattle (lean)"
$ws.Range("B5").Value = 221.7
$ws.Range("C5").Value = 450.0

$ws.Range("A6").Value = "Cattle (fat)"
$ws.Range("B6").Value = 36.8
$ws.Range("C6").Value = 250.0

$ws.Range("A7").Value = "Cattle (organs)"
$ws.Range("B7").Value = 214.6
$ws.Range("C7").Value = 300.0

$ws.Range("A8").Value = "Milk (whole)"
$ws.Range("B8").Value = 327.9
$ws.Range("C8").Value = 200.0

$ws.Range("A9").Value = "Sugar (beets)"
$ws.Range("B9").Value = 64.9
$ws.Range("C9").Value = 250.0

# ---------------------------------------------------------------------------
# 4) Sheet 3 ("II - OS")
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = 1367.0
$ws.Range("C2").Value = 2100.0

$ws.Range("A3").Value = "Potatoes"
$ws.Range("B3").Value = 574.7
$ws.Range("C3").Value = 500.0

$ws.Range("A4").Value = "Wheat (hard red spring)"
$ws.Range("B4").Value = 85.1
$ws.Range("C4").Value = 280.0

$ws.Range("A5").Value = "Barley (pearled)"
$ws.Range("B5").Value = 24.4
$ws.Range("C5").Value = 30.0

$ws.Range("A6").Value = "Canola oil"
$ws.Range("B6").Value = 20.4
$ws.Range("C6").Value = 180.0

$ws.Range("A7").Value = "Rice (brown)"
$ws.Range("B7").Value = 81.3
$ws.Range("C7").Value = 100.0

$ws.Range("A8").Value = "Soy flour"
$ws.Range("B8").Value = 51.8
$ws.Range("C8").Value = 225.0

$ws.Range("A9").Value = "Soybeans"
$ws.Range("B9").Value = 29.1
$ws.Range("C9").Value = 50.0

$ws.Range("A10").Value = "Anchovy (raw)"
$ws.Range("B10").Value = 76.3
$ws.Range("C10").Value = 100.0

$ws.Range("A11").Value = "Cattle (lean)"
$ws.Range("B11").Value = 73.9
$ws.Range("C11").Value = 150.0

$ws.Range("A12").Value = "Cattle (fat)"
$ws.Range("B12").Value = 16.2
$ws.Range("C12").Value = 110.0

$ws.Range("A13").Value = "Cattle (organs)"
$ws.Range("B13").Value = 57.2
$ws.Range("C13").Value = 80.0

$ws.Range("A14").Value = "Milk (whole)"
$ws.Range("B14").Value = 204.9
$ws.Range("C14").Value = 125.0

$ws.Range("A15").Value = "Emi-tsunomata (dry)"
$ws.Range("B15").Value = 38.6
$ws.Range("C15").Value = 100.0

$ws.Range("A16").Value = "Laver (dry)"
$ws.Range("B16").Value = 33.3
$ws.Range("C16").Value = 70.0

# ---------------------------------------------------------------------------
# 5) Sheet 4 ("II - AS")
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = 1338.0
$ws.Range("C2").Value = 2100.0

$ws.Range("A3").Value = "Potatoes"
$ws.Range("B3").Value = 574.7
$ws.Range("C3").Value = 500.0

$ws.Range("A4").Value = "Wheat (hard red spring)"
$ws.Range("B4").Value = 85.1
$ws.Range("C4").Value = 280.0

$ws.Range("A5").Value = "Barley (pearled)"
$ws.Range("B5").Value = 24.4
$ws.Range("C5").Value = 30.0

$ws.Range("A6").Value = "Canola oil"
$ws.Range("B6").Value = 20.4
$ws.Range("C6").Value = 180.0

$ws.Range("A7").Value = "Rice (brown)"
$ws.Range("B7").Value = 81.3
$ws.Range("C7").Value = 100.0

$ws.Range("A8").Value = "Soy flour"
$ws.Range("B8").Value = 51.8
$ws.Range("C8").Value = 225.0

$ws.Range("A9").Value = "Soybeans"
$ws.Range("B9").Value = 29.1
$ws.Range("C9").Value = 50.0

$ws.Range("A10").Value = "Anchovy (raw)"
$ws.Range("B10").Value = 76.3
$ws.Range("C10").Value = 100.0

$ws.Range("A11").Value = "Cattle (lean)"
$ws.Range("B11").Value = 73.9
$ws.Range("C11").Value = 150.0

$ws.Range("A12").Value = "Cattle (fat)"
$ws.Range("B12").Value = 16.2
$ws.Range("C12").Value = 110.0

$ws.Range("A13").Value = "Cattle (organs)"
$ws.Range("B13").Value = 57.2
$ws.Range("C13").Value = 80.0

$ws.Range("A14").Value = "Milk (whole)"
$ws.Range("B14").Value = 204.9
$ws.Range("C14").Value = 125.0

$ws.Range("A15").Value = "Lignocellulosic sugar"
$ws.Range("B15").Value = 42.5
$ws.Range("C15").Value = 170.0

# ---------------------------------------------------------------------------
# 6) Sheet 5 ("III- OS")
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(5)
$ws.Range("B2").Value = 1232.0
$ws.Range("C2").Value = 2100.0

$ws.Range("A3").Value = "Potatoes"
$ws.Range("B3").Value = 574.7
$ws.Range("C3").Value = 500.0

$ws.Range("A4").Value = "Wheat (hard red spring)"
$ws.Range("B4").Value = 76.0
$ws.Range("C4").Value = 250.0

$ws.Range("A5").Value = "Barley (pearled)"
$ws.Range("B5").Value = 40.7
$ws.Range("C5").Value = 50.0

$ws.Range("A6").Value = "Canola oil"
$ws.Range("B6").Value = 22.6
$ws.Range("C6").Value = 200.0

$ws.Range("A7").Value = "Rice (brown)"
$ws.Range("B7").Value = 81.3
$ws.Range("C7").Value = 100.0

$ws.Range("A8").Value = "Corn flour (whole-grain)"
$ws.Range("B8").Value = 41.6
$ws.Range("C8").Value = 150.0

$ws.Range("A9").Value = "Corn"
$ws.Range("B9").Value = 140.6
$ws.Range("C9").Value = 135.0

$ws.Range("A10").Value = "Soy flour"
$ws.Range("B10").Value = 46.1
$ws.Range("C10").Value = 200.0

$ws.Range("A11").Value = "Soybeans"
$ws.Range("B11").Value = 46.5
$ws.Range("C11").Value = 80.0

$ws.Range("A12").Value = "Milk (whole)"
$ws.Range("B12").Value = 41.0
$ws.Range("C12").Value = 25.0

$ws.Range("A13").Value = "Spirulina (dry)"
$ws.Range("B13").Value = 17.2
$ws.Range("C13").Value = 50.0

$ws.Range("A14").Value = "Emi-tsunomata (dry)"
$ws.Range("B14").Value = 38.6
$ws.Range("C14").Value = 100.0

$ws.Range("A15").Value = "Laver (dry)"
$ws.Range("B15").Value = 21.4
$ws.Range("C15").Value = 45.0

$ws.Range("A16").Value = "Wakame (dry)"
$ws.Range("B16").Value = 5.6
$ws.Range("C16").Value = 15.0

$ws.Range("A17").Value = "Bacteria (methane)"
$ws.Range("B17").Value = 38.0
$ws.Range("C17").Value = 200.0

# ---------------------------------------------------------------------------
# 7) Sheet 6 ("III - AS")
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(6)
$ws.Range("B2").Value = 1121.0
$ws.Range("C2").Value = 2100.0

$ws.Range("A3").Value = "Potatoes"
$ws.Range("B3").Value = 574.7
$ws.Range("C3").Value = 500.0

$ws.Range("A4").Value = "Wheat (hard red spring)"
$ws.Range("B4").Value = 106.4
$ws.Range("C4").Value = 350.0

$ws.Range("A5").Value = "Barley (pearled)"
$ws.Range("B5").Value = 40.7
$ws.Range("C5").Value = 50.0

$ws.Range("A6").Value = "Canola oil"
$ws.Range("B6").Value = 28.3
$ws.Range("C6").Value = 250.0

$ws.Range("A7").Value = "Cattle (fat)"
$ws.Range("B7").Value = 29.4
$ws.Range("C7").Value = 200.0

$ws.Range("A8").Value = "Milk (whole)"
$ws.Range("B8").Value = 163.9
$ws.Range("C8").Value = 100.0

$ws.Range("A9").Value = "Spirulina (dry)"
$ws.Range("B9").Value = 17.2
$ws.Range("C9").Value = 50.0

$ws.Range("A10").Value = "Emi-tsunomata (dry)"
$ws.Range("B10").Value = 38.6
$ws.Range("C10").Value = 100.0

$ws.Range("A11").Value = "Laver (dry)"
$ws.Range("B11").Value = 23.8
$ws.Range("C11").Value = 50.0

$ws.Range("A12").Value = "Wakame (dry)"
$ws.Range("B12").Value = 9.3
$ws.Range("C12").Value = 25.0

$ws.Range("A13").Value = "Bacteria (methane)"
$ws.Range("B13").Value = 57.0
$ws.Range("C13").Value = 300.0

$ws.Range("A14").Value = "Lignocellulosic sugar"
$ws.Range("B14").Value = 31.3
$ws.Range("C14").Value = 125.0

# Page setup additions on sheet 6 (fit-to-page, print options, margins)
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.PrintGridlines = $true
$ws.PageSetup.CenterHorizontally = $true
$ws.PageSetup.BottomMargin = 0.75
$ws.PageSetup.HeaderMargin = 0.0
$ws.PageSetup.FooterMargin = 0.0
$ws.PageSetup.LeftMargin = 0.7
$ws.PageSetup.RightMargin = 0.7
$ws.PageSetup.TopMargin = 0.75
$ws.PageSetup.FitToPagesTall = 0
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
$ws.PageSetup.Order = 1
